# "Update 2021 HWL2 First Batch"
#
# Adds year columns 2016-2050 (35 new years) to the Czechoslovakia labourers'
# real wage workbook, and tidies up the wording of the text citation.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Data Clio Infra Format" sheet - wide/Clio format.
#    Existing year columns run 1500 (col I) .. 2015 (col TD). Append 35 more
#    year columns (2016 .. 2050) right after TD, i.e. TE .. UM. Row 2 (the
#    only data row) is left blank for all of the new columns.
# ---------------------------------------------------------------------------
$wsWide = $wb.Worksheets.Item("Data Clio Infra Format")

$firstNewCol = 525   # column TE
for ($i = 0; $i -lt 35; $i++) {
    $year = 2016 + $i
    $col = $firstNewCol + $i
    # Leading apostrophe forces text ("string") storage, matching the other
    # year-header cells (1500 .. 2015) which are stored as text, not numbers.
    $wsWide.Cells.Item(1, $col).Value = "'" + $year.ToString()
}

# ---------------------------------------------------------------------------
# 2) "Data Long Format" sheet - long format.
#    Columns are: A country name, B Borders Start Year, C Borders End Year,
#    D Indicator, E year, F value.
#    Insert 35 blank columns before the old E ("year") column and label them
#    with the same new years (2016 .. 2050); this pushes "year"/"value" out
#    to AN/AO.
# ---------------------------------------------------------------------------
$wsLong = $wb.Worksheets.Item("Data Long Format")

$wsLong.Range("E1:AM1").EntireColumn.Insert()

$firstNewColLong = 5   # column E
for ($i = 0; $i -lt 35; $i++) {
    $year = 2016 + $i
    $col = $firstNewColLong + $i
    $wsLong.Cells.Item(1, $col).Value = "'" + $year.ToString()
}

# ---------------------------------------------------------------------------
# 3) "Metadata" sheet - fix the text citation punctuation (Oxford comma
#    before "and").
# ---------------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("C3").Value = "Zwart, Pim de, Bas van Leeuwen, and Jieli van Leeuwen-Li (2015). Labourers Real Wage. http://hdl.handle.net/10622/QK8VRF, accessed via the Clio Infra website."
